$wb = $excel.ActiveWorkbook

# The automation used to write weather data starting at A1 of a brand new
# sheet (previously it only wrote the "Location" column into the existing
# sheet). Add the new sheet that receives the full write range and place it
# before the original sheet, matching the new workbook tab order.
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Arkusz2"

# Values coming from the weather API are written as text (quantities like
# "10", "66" keep their shared-string/text type instead of becoming numbers).
$newSheet.Range("A1:E6").NumberFormat = "@"

$data = @(
    @("Location", "Temperature", "Sky", "Humidity", "Wind"),
    @("london", "10", "Klart", "66", "7 m/s Sydvest"),
    @("warsaw", "13", "Skyet", "100", "3 m/s Syd"),
    @("paris", "11", "Klart", "76", "3 m/s Vest"),
    @("madrit", "15", "Klart", "88", "1 m/s Nordvest"),
    @("budapest", "12", "Mest skyet", "94", "2 m/s Øst")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    for ($c = 0; $c -lt $data[$r].Length; $c++) {
        $newSheet.Cells.Item($r + 1, $c + 1).Value = $data[$r][$c]
    }
}

# The original sheet keeps holding just the location names (column A) plus
# its header row, but now shares the same cell style as the new sheet.
$oldSheet = $wb.Worksheets.Item("Arkusz1")
$oldSheet.Range("A1:E1").NumberFormat = "@"
$oldSheet.Range("A2:A6").NumberFormat = "@"

# Restore each sheet's last selection / active cell.
$oldSheet.Range("A2").Select()
$newSheet.Activate()
$newSheet.Range("I9").Select()
